# Revision: add "chemical_recycling_pyrolysis" parameter row and shift
# the following rows down (commit: "revision, added pyrolysis and
# additional figures").
#
# In the user-inputs parameter table on Sheet1, a new boolean parameter
# "chemical_recycling_pyrolysis" (set to TRUE) is inserted directly below
# the existing "chemical_recycling_gasification" row (row 9), pushing all
# subsequent parameter rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 10 (shifts rows 10:24 down to 11:25).
$ws.Rows("10:10").Insert()

# Populate the new row with the new parameter.
$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
